$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 562.25
$ws.Range("I32").Value = 300
$ws.Range("J32").Value = 599.7143
$ws.Range("K32").Value = 300
$ws.Range("L32").Value = 599.7143
$ws.Range("M32").Value = 26
$ws.Range("N32").Value = -1251.7143
$ws.Range("H76").Value = 7621.9287
$ws.Range("I76").Value = 8652.947
$ws.Range("K76").Value = 8652.947
$ws.Range("M76").Value = -8337.947
$ws.Range("H79").Value = 7621.9287
$ws.Range("I79").Value = 8652.947
$ws.Range("K79").Value = 8652.947
$ws.Range("M79").Value = -7560.947
$ws.Range("H99").Value = 906.4
$ws.Range("I99").Value = 685.4286
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 2056.2858
$ws.Range("L99").Value = 12000
$ws.Range("M99").Value = -558.2857999999997
$ws.Range("N99").Value = -14996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 21000
$ws.Range("J44").Value = 21000
$ws.Range("L44").Value = 21000
$ws.Range("N44").Value = -21976
$ws.Range("H122").Value = 1974255.1
$ws.Range("I122").Value = 1974255.1
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5922765.300000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5920315.300000001
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 5265.3335
$ws.Range("I54").Value = 2125
$ws.Range("J54").Value = 11546
$ws.Range("K54").Value = 2125
$ws.Range("L54").Value = 11546
$ws.Range("M54").Value = -1641
$ws.Range("N54").Value = -12514
$ws.Range("H64").Value = 231
$ws.Range("I64").Value = 164.625
$ws.Range("J64").Value = 284.1
$ws.Range("K64").Value = 164.625
$ws.Range("L64").Value = 284.1
$ws.Range("M64").Value = 60.375
$ws.Range("N64").Value = -734.1
$ws.Range("H67").Value = 231
$ws.Range("I67").Value = 164.625
$ws.Range("J67").Value = 284.1
$ws.Range("K67").Value = 164.625
$ws.Range("L67").Value = 284.1
$ws.Range("M67").Value = 615.375
$ws.Range("N67").Value = -1844.1
$ws.Range("H86").Value = 1661.8108
$ws.Range("I86").Value = 1615.8438
$ws.Range("J86").Value = 1956
$ws.Range("K86").Value = 1615.8438
$ws.Range("L86").Value = 1956
$ws.Range("M86").Value = -492.8438000000001
$ws.Range("N86").Value = -4202
$ws.Range("H89").Value = 1661.8108
$ws.Range("I89").Value = 1615.8438
$ws.Range("J89").Value = 1956
$ws.Range("K89").Value = 8079.219000000001
$ws.Range("L89").Value = 9780
$ws.Range("M89").Value = -2463.219000000001
$ws.Range("N89").Value = -21012

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17564.75
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 17564.75
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 17564.75
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -18154.75
$ws.Range("H34").Value = 17564.75
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 17564.75
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 17564.75
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -17968.75
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H134").Value = 4256.857
$ws.Range("I134").Value = 4577.294
$ws.Range("J134").Value = 2895
$ws.Range("K134").Value = 13731.882
$ws.Range("L134").Value = 8685
$ws.Range("M134").Value = -11196.882
$ws.Range("N134").Value = -13755

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 218343.12
$ws.Range("I5").Value = 409.33334
$ws.Range("J5").Value = 323794.97
$ws.Range("K5").Value = 1228.00002
$ws.Range("L5").Value = 971384.9099999999
$ws.Range("M5").Value = -1116.00002
$ws.Range("N5").Value = -971608.9099999999
$ws.Range("H87").Value = 4750
$ws.Range("I87").Value = 4750
$ws.Range("K87").Value = 14250
$ws.Range("M87").Value = -13002
$ws.Range("H90").Value = 4750
$ws.Range("I90").Value = 4750
$ws.Range("K90").Value = 42750
$ws.Range("M90").Value = -36510
$ws.Range("H131").Value = 961.62744
$ws.Range("J131").Value = 1019.4186
$ws.Range("L131").Value = 3058.2558
$ws.Range("N131").Value = -13138.2558
$ws.Range("H135").Value = 218343.12
$ws.Range("I135").Value = 409.33334
$ws.Range("J135").Value = 323794.97
$ws.Range("K135").Value = 3684.00006
$ws.Range("L135").Value = 2914154.73
$ws.Range("M135").Value = -1149.00006
$ws.Range("N135").Value = -2919224.73

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5520.4
$ws.Range("I70").Value = 5961.3887
$ws.Range("K70").Value = 5961.3887
$ws.Range("M70").Value = -5691.3887
$ws.Range("H73").Value = 5520.4
$ws.Range("I73").Value = 5961.3887
$ws.Range("K73").Value = 5961.3887
$ws.Range("M73").Value = -5025.3887
$ws.Range("H122").Value = 2495485.8
$ws.Range("I122").Value = 4631803
$ws.Range("J122").Value = 3115.5
$ws.Range("K122").Value = 13895409
$ws.Range("L122").Value = 9346.5
$ws.Range("M122").Value = -13892959
$ws.Range("N122").Value = -14246.5
$ws.Range("H123").Value = 21806.654
$ws.Range("J123").Value = 21806.654
$ws.Range("L123").Value = 21806.654
$ws.Range("N123").Value = -26706.654

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1057.1515
$ws.Range("I22").Value = 150.33333
$ws.Range("J22").Value = 1575.3334
$ws.Range("K22").Value = 150.33333
$ws.Range("L22").Value = 1575.3334
$ws.Range("M22").Value = 144.66667
$ws.Range("N22").Value = -2165.3334
$ws.Range("H27").Value = 1057.1515
$ws.Range("I27").Value = 150.33333
$ws.Range("J27").Value = 1575.3334
$ws.Range("K27").Value = 150.33333
$ws.Range("L27").Value = 1575.3334
$ws.Range("M27").Value = -43.33332999999999
$ws.Range("N27").Value = -1789.3334
$ws.Range("H82").Value = 395526.94
$ws.Range("J82").Value = 58930.832
$ws.Range("L82").Value = 58930.832
$ws.Range("N82").Value = -59652.832
$ws.Range("H85").Value = 395526.94
$ws.Range("J85").Value = 58930.832
$ws.Range("L85").Value = 58930.832
$ws.Range("N85").Value = -61426.832
$ws.Range("H98").Value = 28838.75
$ws.Range("J98").Value = 28838.75
$ws.Range("L98").Value = 28838.75
$ws.Range("N98").Value = -34828.75
$ws.Range("H122").Value = 2810132.2
$ws.Range("I122").Value = 3403401.5
$ws.Range("J122").Value = 1252800.6
$ws.Range("K122").Value = 10210204.5
$ws.Range("L122").Value = 3758401.8
$ws.Range("M122").Value = -10207754.5
$ws.Range("N122").Value = -3763301.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1098.1111
$ws.Range("I122").Value = 1125.5
$ws.Range("J122").Value = 879
$ws.Range("K122").Value = 3376.5
$ws.Range("L122").Value = 2637
$ws.Range("M122").Value = -926.5
$ws.Range("N122").Value = -7537
$ws.Range("H123").Value = 29350
$ws.Range("J123").Value = 29350
$ws.Range("L123").Value = 29350
$ws.Range("N123").Value = -39150
